# "Spelling error corrected & style simplified"
#
# 1. Rename the "table" sheet to "Table" and add two new blank sheets
#    (Sheet2, Sheet3) right after it.
# 2. Fix the "clowdy" -> "cloudy" typo in D8.
# 3. Simplify the formatting: every coloured/bold font is reset to the
#    theme's automatic text colour (theme 1) and every coloured fill is
#    removed, while keeping the bold/italic/font-family, number format and
#    alignment that was already in place for each block of cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. sheet name + extra sheets -------------------------------------------------
$ws.Name = "Table"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sheet2"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Sheet3"

# --- 2. fix the typo --------------------------------------------------------------
$ws.Range("D8").Formula = '="cloudy"'

# --- 3. style simplification -------------------------------------------------------

# Row 1 (A1:E1): keep the bold Courier New / numeric "0" format / centered
# look, just drop the per-column accent colour back to automatic (theme 1).
$headerRow = $ws.Range("A1:E1")
$headerRow.Font.ThemeColor = 1

# Rows 2-5 (A:D): plain text cells with a coloured fill per column -> drop
# the fill and keep left/center text alignment.
$plainBlock = $ws.Range("A2:D5")
$plainBlock.ClearFormats()
$plainBlock.NumberFormat = "@"
$plainBlock.HorizontalAlignment = -4131
$plainBlock.VerticalAlignment = -4108

# Row 6 (A6:D6): bold italic Arial header cells -> drop fill & accent colour.
$row6Left = $ws.Range("A6:D6")
$row6Left.ClearFormats()
$row6Left.NumberFormat = "@"
$row6Left.Font.Bold = $true
$row6Left.Font.Italic = $true
$row6Left.Font.Name = "Arial"
$row6Left.Font.ThemeColor = 1
$row6Left.HorizontalAlignment = -4131
$row6Left.VerticalAlignment = -4108

# Row 6, column E: same font treatment but right aligned.
$row6Right = $ws.Range("E6")
$row6Right.ClearFormats()
$row6Right.NumberFormat = "@"
$row6Right.Font.Bold = $true
$row6Right.Font.Italic = $true
$row6Right.Font.Name = "Arial"
$row6Right.Font.ThemeColor = 1
$row6Right.HorizontalAlignment = -4152
$row6Right.VerticalAlignment = -4108

# Rows 7-8 (A7:D7, B8:D8): bold (non-italic) Arial -> drop fill & accent colour.
# (Kept as two separate statements per row: multi-area Range property
# assignment only reliably touches the first area.)
$row7Left = $ws.Range("A7:D7")
$row7Left.ClearFormats()
$row7Left.NumberFormat = "@"
$row7Left.Font.Bold = $true
$row7Left.Font.Name = "Arial"
$row7Left.Font.ThemeColor = 1
$row7Left.HorizontalAlignment = -4131
$row7Left.VerticalAlignment = -4108

$row8Left = $ws.Range("B8:D8")
$row8Left.ClearFormats()
$row8Left.NumberFormat = "@"
$row8Left.Font.Bold = $true
$row8Left.Font.Name = "Arial"
$row8Left.Font.ThemeColor = 1
$row8Left.HorizontalAlignment = -4131
$row8Left.VerticalAlignment = -4108

# Column E, rows 7-9: same bold Arial font but right aligned.
$colERight = $ws.Range("E7:E9")
$colERight.ClearFormats()
$colERight.NumberFormat = "@"
$colERight.Font.Bold = $true
$colERight.Font.Name = "Arial"
$colERight.Font.ThemeColor = 1
$colERight.HorizontalAlignment = -4152
$colERight.VerticalAlignment = -4108

$ws.Activate()
